$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append rows 302..328, continuing the daily date series in column A
# (serial 44376 .. 44402) with zeros in columns B, C, D, matching the
# style/format already used for the preceding rows.

$startRow = 302
$startSerial = 44376
$endSerial = 44402

$templateRow = $startRow - 1

$row = $startRow
for ($serial = $startSerial; $serial -le $endSerial; $serial++) {
    # Copy the previous row's A:D formatting (incl. the date cell's
    # existing style) down before writing the new values, so no new
    # cellXf entries get minted in styles.xml and no stray formatted
    # cells appear outside the used range.
    $srcRange = $ws.Range("A" + $templateRow + ":D" + $templateRow)
    $dstRange = $ws.Range("A" + $row + ":D" + $row)
    $srcRange.Copy($dstRange)

    $ws.Cells.Item($row, 1).Value = $serial
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
    $row++
}
